# Edit script: "Not respecting bolds and whitespaces"
#
# 1) Resize / relayout the contract-summary table (7 columns, 2 rows):
#    new column widths, fixed layout, narrower overall table width.
# 2) Fix the wording "... asciende a un total de #...#" -> "... asciende a  #...#"
#    (drop "un total de", matching bold/whitespace fix from the commit).
# 3) Clean up the "#Telefono#" placeholder in the header so it is a single,
#    consistently formatted run (no stray spell-check markers splitting it up).
# 4) Drop the stale "_GoBack" bookmark left over from the previous edit session.

$d = $word.ActiveDocument

# --- 1) Table: widths, table width, fixed layout ---------------------------
$t = $d.Tables(1)

$t.AllowAutoFit = $false

$t.Columns(1).Width = 67.05
$t.Columns(2).Width = 74.75
$t.Columns(3).Width = 85.05
$t.Columns(4).Width = 56.7
$t.Columns(5).Width = 99.2
$t.Columns(6).Width = 70.9
$t.Columns(7).Width = 70.85

$t.PreferredWidthType = 3
$t.PreferredWidth = 524.5

# --- 2) Wording fix: remove "un total de" -----------------------------------
$d.Content.Find.Execute(" asciende a un total de ", $true, $false, $false, `
    $false, $false, $true, 1, $false, " asciende a  ", 2)

# --- 3) Header: merge "#" + "Telefono" + "#" into a single run -------------
$hdr = $d.Sections(1).Headers(2).Range
$hdr.Find.Execute("#Telefono#", $false, $false, $false, $false, $false, `
    $true, 1, $false, "#Telefono#", 2)

# --- 4) Remove the leftover _GoBack bookmark --------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
